$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 450, shifting the existing rows 450-548
# down to 452-550 (matching the diff, which shows a new week's worth of
# data - two new records dated 44995 - prepended to this block, with the
# former rows sliding down by two positions).
$ws.Rows.Item(450).EntireRow.Insert()
$ws.Rows.Item(450).EntireRow.Insert()

# Populate the first new row (450)
$ws.Range("A450").Value = 10
$ws.Range("B450").Value = "Vega Modelo de Temuco"
$ws.Range("C450").Value = "La Araucanía"
$ws.Range("D450").Value = 44995
$ws.Range("E450").Value = 9
$ws.Range("F450").Value = 100114014
$ws.Range("G450").Value = "Betarraga"
$ws.Range("H450").Value = "Sin especificar"
$ws.Range("I450").Value = "Primera"
$ws.Range("J450").Value = 110
$ws.Range("K450").Value = 10000
$ws.Range("L450").Value = 10000
$ws.Range("M450").Value = 10000
$ws.Range("N450").Value = "$/docena de paquetes"
$ws.Range("O450").Value = "Provincia de Cautín"
$ws.Range("P450").Value = 833
$ws.Range("Q450").Value = 12
$ws.Range("R450").Value = "Hortaliza"

# Populate the second new row (451)
$ws.Range("A451").Value = 10
$ws.Range("B451").Value = "Vega Modelo de Temuco"
$ws.Range("C451").Value = "La Araucanía"
$ws.Range("D451").Value = 44995
$ws.Range("E451").Value = 9
$ws.Range("F451").Value = 100114014
$ws.Range("G451").Value = "Betarraga"
$ws.Range("H451").Value = "Sin especificar"
$ws.Range("I451").Value = "Primera"
$ws.Range("J451").Value = 75
$ws.Range("K451").Value = 10000
$ws.Range("L451").Value = 10000
$ws.Range("M451").Value = 10000
$ws.Range("N451").Value = "$/docena de paquetes"
$ws.Range("O451").Value = "Región del Maule"
$ws.Range("P451").Value = 833
$ws.Range("Q451").Value = 12
$ws.Range("R451").Value = "Hortaliza"
